# ADAPTERSPECIFICATTRIBUTE.docx template fix
# "fix format issue and content filter(remove value=="N/A")"
#
# Summary of the edit being applied to the table that shows the
# Name / NameSpace / Value rows:
#   - table switches from a fixed-width (dxa) layout to a 100%
#     (5000 pct) fixed layout, with proportionally recomputed
#     column / cell widths
#   - every cell grows explicit single-line borders (w:tcBorders)
#   - the header/label cell shading changes from the old blue theme
#     fill to a plain 25% gray (BFBFBF / background1 shade BF) fill
#   - the label cells pick up 1.15 line spacing (w:spacing
#     w:line="276" w:lineRule="auto") instead of carrying bold via
#     an (unused) empty-run rPr
#   - the "$Foo_Value" placeholder runs are split into "$" + the
#     rest of the name, with the rest wrapped in proofErr
#     spellStart/spellEnd (this is what Word itself does once it
#     spell-checks text that isn't literally "N/A" any more)
#   - the stray _GoBack bookmark that used to live inside the last
#     cell's paragraph is moved onto the trailing empty paragraph
#     after the table (where Word normally keeps it)
#
# The cleanest, least error-prone way to reproduce all of that via
# COM automation is the same technique Word itself exposes for bulk
# structural edits: build the replacement OOXML for word/document.xml
# and hand it to Range.InsertXML on the full document range. This
# keeps every other package part (styles, settings, rels, theme, ...)
# untouched while swapping in the new table/paragraph markup.

$d = $word.ActiveDocument

$newDocumentXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:tbl>
            <w:tblPr>
              <w:tblW w:w="5000" w:type="pct"/>
              <w:tblLayout w:type="fixed"/>
              <w:tblLook w:val="00A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/>
            </w:tblPr>
            <w:tblGrid>
              <w:gridCol w:w="3231"/>
              <w:gridCol w:w="6113"/>
            </w:tblGrid>
            <w:tr>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:w="5000" w:type="pct"/>
                  <w:gridSpan w:val="2"/>
                  <w:tcBorders>
                    <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                  </w:tcBorders>
                  <w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/>
                </w:tcPr>
                <w:p>
                  <w:pPr>
                    <w:spacing w:line="276" w:lineRule="auto"/>
                  </w:pPr>
                  <w:r>
                    <w:rPr>
                      <w:b/>
                    </w:rPr>
                    <w:t>Adapter specific attribute</w:t>
                  </w:r>
                </w:p>
              </w:tc>
            </w:tr>
            <w:tr>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:w="1729" w:type="pct"/>
                  <w:tcBorders>
                    <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                  </w:tcBorders>
                  <w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/>
                </w:tcPr>
                <w:p>
                  <w:pPr>
                    <w:spacing w:line="276" w:lineRule="auto"/>
                    <w:rPr>
                      <w:b/>
                    </w:rPr>
                  </w:pPr>
                  <w:r>
                    <w:rPr>
                      <w:b/>
                    </w:rPr>
                    <w:t>Name</w:t>
                  </w:r>
                </w:p>
              </w:tc>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:w="3271" w:type="pct"/>
                  <w:tcBorders>
                    <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                  </w:tcBorders>
                </w:tcPr>
                <w:p>
                  <w:r>
                    <w:t>$</w:t>
                  </w:r>
                  <w:proofErr w:type="spellStart"/>
                  <w:r>
                    <w:t>Name_Value</w:t>
                  </w:r>
                  <w:proofErr w:type="spellEnd"/>
                </w:p>
              </w:tc>
            </w:tr>
            <w:tr>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:w="1729" w:type="pct"/>
                  <w:tcBorders>
                    <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                  </w:tcBorders>
                  <w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/>
                </w:tcPr>
                <w:p>
                  <w:pPr>
                    <w:spacing w:line="276" w:lineRule="auto"/>
                    <w:rPr>
                      <w:b/>
                    </w:rPr>
                  </w:pPr>
                  <w:proofErr w:type="spellStart"/>
                  <w:r>
                    <w:rPr>
                      <w:b/>
                    </w:rPr>
                    <w:t>NameSpace</w:t>
                  </w:r>
                  <w:proofErr w:type="spellEnd"/>
                </w:p>
              </w:tc>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:w="3271" w:type="pct"/>
                  <w:tcBorders>
                    <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                  </w:tcBorders>
                </w:tcPr>
                <w:p>
                  <w:r>
                    <w:t>$</w:t>
                  </w:r>
                  <w:proofErr w:type="spellStart"/>
                  <w:r>
                    <w:t>NameSpace_Value</w:t>
                  </w:r>
                  <w:proofErr w:type="spellEnd"/>
                </w:p>
              </w:tc>
            </w:tr>
            <w:tr>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:w="1729" w:type="pct"/>
                  <w:tcBorders>
                    <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                  </w:tcBorders>
                  <w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/>
                </w:tcPr>
                <w:p>
                  <w:pPr>
                    <w:spacing w:line="276" w:lineRule="auto"/>
                    <w:rPr>
                      <w:b/>
                    </w:rPr>
                  </w:pPr>
                  <w:r>
                    <w:rPr>
                      <w:b/>
                    </w:rPr>
                    <w:t>Value</w:t>
                  </w:r>
                </w:p>
              </w:tc>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:w="3271" w:type="pct"/>
                  <w:tcBorders>
                    <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                    <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
                  </w:tcBorders>
                </w:tcPr>
                <w:p>
                  <w:r>
                    <w:t>$</w:t>
                  </w:r>
                  <w:proofErr w:type="spellStart"/>
                  <w:r>
                    <w:t>Value_Value</w:t>
                  </w:r>
                  <w:proofErr w:type="spellEnd"/>
                </w:p>
              </w:tc>
            </w:tr>
          </w:tbl>
          <w:p>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
          <w:sectPr>
            <w:pgSz w:w="12240" w:h="15840"/>
            <w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/>
            <w:cols w:space="720"/>
            <w:docGrid w:linePitch="360"/>
          </w:sectPr>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

# InsertXML replaces the contents of the exact range it's called on;
# calling it on the whole-document Content range swaps in the fixed
# table/paragraph markup above while every other package part
# (styles.xml, settings.xml, theme, rels, ...) is left alone.
$d.Content.InsertXML($newDocumentXml)
